$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title at the top of the document.
# ---------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Insert a new paragraph ("Play Age of The Gods Prince of Olympus
#    Free | Slot Review", bold) right before the final "Prompt: ..."
#    paragraph at the end of the document.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPar = $d.Paragraphs($count)
$lastStart = $lastPar.Range.Start
$insertionPoint = $d.Range($lastStart, $lastStart)
$insertionPoint.InsertAfter("Play Age of The Gods Prince of Olympus Free | Slot Review" + [char]13)

$count = $d.Paragraphs.Count
$newPar = $d.Paragraphs($count - 1)
$newParRange = $newPar.Range
$newTextRange = $d.Range($newParRange.Start, $newParRange.End - 1)
$newTextRange.Font.Bold = 1

# ---------------------------------------------------------------------
# 3) Replace the old AI-image "Prompt: ..." text (now the very last
#    paragraph) with the new meta-description copy, keeping the
#    existing italic run formatting untouched.
# ---------------------------------------------------------------------
$oldText = 'Prompt: Create a feature image for "Age of The Gods Prince of Olympus" that captures the essence of the game - the powerful son of Zeus - Hercules, and the fun and excitement that this slot game brings. Specifications: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses - The image should be eye-catching and playful - The image should include the game title "Age of The Gods Prince of Olympus" Keep in mind the theme of the game, which revolves around Greek mythology and Hercules. The Maya warrior with glasses adds an element of surprise and fun to the image, resulting in a unique and captivating piece of art. The image should aim to draw players in and encourage them to try out the game.'
$newText = 'Discover the world of the Greek hero Hercules and try your luck with exciting bonus rounds and Progressive Jackpots in Age of The Gods Prince of Olympus. Play for free now!'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
